$d = $word.ActiveDocument

# Build the replacement runs (with correct per-run formatting, including w:lang on some)
# as a scratch paragraph inserted via raw OOXML, then lift its FormattedText onto the
# target range, then remove the scratch paragraph.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>veoma</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> brzo</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>. Programeri su “u hodu” otkrivali kako najbolje da primene postoje</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">će šablone u programiranju na novi kontekst i uporedo sa tim su se razvijali alati i biblioteke koji pojednostavljuju i ubrzavaju programiranje Android aplikacija. Vremenom </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">je pređen </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>put od aplikacija čija je sva logika sadržana u aktivnostima aplikacije (eng. Activity) preko korišćenja Model-pogled-kontroler šablona (eng. Model-View-Controller MVC) sve do danas aktuelnih šablona Model-kontroler-prezenter (eng. Model-View-Presenter MVP)  i Model-Pogled-Pogled-Model (eng. Model-View-View-Model MVVM).</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve"> Takođe zbog sve složenijih projekata a prateći SOLID</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> princip inverzije kontrole zapo</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">četo je korišćenje bilblioteka koje omogućavaju </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>umetanje zavisnosti</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">kao što su RoboGuice i Dagger. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$end = $d.Content.End
$scratchAnchor = $d.Range($end, $end)
$scratchAnchor.InsertXML($xml)

$scratchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratchTextRange = $d.Range($scratchPara.Range.Start, $scratchPara.Range.End - 1)
Write-Host "scratch text: [$($scratchTextRange.Text)]"
$ft = $scratchTextRange.FormattedText

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("vrlo brzo")
Write-Host "found vrlo brzo: $found"
$target = $find.Parent
$target.FormattedText = $ft
Write-Host "replaced target text"

$scratchPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$delRange = $d.Range($scratchPara2.Range.Start - 1, $scratchPara2.Range.End)
$delRange.Delete()
Write-Host "cleaned up scratch paragraph"
